# "Removed the individual component slides."
#
# The deck had five "Component: ..." slides (pubsub, Broker, Adapter,
# client endpoint, Service endpoint) sitting right before the four
# "Deployment scenarios" slides. This commit drops those five component
# slides, so the deployment-scenario slides shift up to take their place.
# It also trims the "Deployment scenarios" overview slide's body text,
# removing the "Major aspect of the "Flexibility" goal" bullet.

$p = $ppt.ActivePresentation

# Slides 14-18 are, in order:
#   14: Component: pubsub
#   15: Component: Broker
#   16: Component: Adapter
#   17: Component: client endpoint
#   18: Component: Service endpoint
# Delete from the back so earlier indices stay valid.
$p.Slides.Item(18).Delete()
$p.Slides.Item(17).Delete()
$p.Slides.Item(16).Delete()
$p.Slides.Item(15).Delete()
$p.Slides.Item(14).Delete()

# What used to be slide 19 ("Deployment scenarios") is now slide 14.
# Remove its first bullet ("Major aspect of the "Flexibility" goal").
$contentShape = $p.Slides.Item(14).Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Delete()
